$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "58.658.85"
$ws.Range("E2").Value = "  +0.73%  "
$ws.Range("D3").Value = "3.155.32"
$ws.Range("E3").Value = "  +0.52%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'529.52"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.08%  "
$ws.Range("D6").Value = "'139.43"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("E8").Value = "  +14.07%  "
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("E10").Value = "  +5.78%  "
$ws.Range("E11").Value = "  +3.62%  "
$ws.Range("D13").Value = "3.698.48"
$ws.Range("E13").Value = "  +0.73%  "
$ws.Range("D14").Value = "'0.0000173"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'25.77"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.79%  "
$ws.Range("D16").Value = "58.720.39"
$ws.Range("E16").Value = "  +0.79%  "
$ws.Range("D17").Value = "'6.26"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.35%  "
$ws.Range("D18").Value = "3.155.67"
$ws.Range("E18").Value = "  +0.42%  "
$ws.Range("D19").Value = "'12.96"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.03%  "
$ws.Range("D20").Value = "'8.13"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.37%  "
$ws.Range("E21").Value = "  +3.47%  "
$ws.Range("E22").Value = "  +1.56%  "
$ws.Range("D23").Value = "'0.997"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.13%  "
$ws.Range("D24").Value = "'0.530"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.75%  "
$ws.Range("E25").Value = "  +0.67%  "
$ws.Range("E26").Value = "  +0.37%  "
$ws.Range("E27").Value = "  -0.26%  "
$ws.Range("D28").Value = "'8.31"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +13.46%  "
$ws.Range("D29").Value = "0.0₃0864"
$ws.Range("E29").Value = "  -1.86%  "
$ws.Range("D30").Value = "'22.26"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.69%  "
$ws.Range("E31").Value = "  -0.17%  "
$ws.Range("E32").Value = "  -1.30%  "
$ws.Range("D33").Value = "'5.12"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.73%  "
$ws.Range("E34").Value = "  +0.53%  "
$ws.Range("E35").Value = "  +3.25%  "
$ws.Range("D36").Value = "'157.99"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.77%  "
$ws.Range("D37").Value = "'1.34"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.40%  "
$ws.Range("D38").Value = "'25.02"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.92%  "
$ws.Range("E39").Value = "  -0.68%  "
$ws.Range("E40").Value = "  +1.82%  "
$ws.Range("D41").Value = "2.620.74"
$ws.Range("E41").Value = "  +4.94%  "
$ws.Range("E42").Value = "  +5.49%  "
$ws.Range("D43").Value = "'0.720"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.56%  "
$ws.Range("D44").Value = "'39.03"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.07%  "
$ws.Range("D45").Value = "'0.0285"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.12%  "
$ws.Range("E46").Value = "  +0.06%  "
$ws.Range("D47").Value = "3.195.28"
$ws.Range("E47").Value = "  +0.60%  "
$ws.Range("E48").Value = "  +12.83%  "
$ws.Range("E49").Value = "  +2.03%  "
$ws.Range("E50").Value = "  -1.80%  "
$ws.Range("D51").Value = "'20.07"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.06%  "
